$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = "2026-02-06 22:30:26"
    3  = "2026-02-06 22:30:28"
    4  = "2026-02-06 22:30:31"
    5  = "2026-02-06 22:30:33"
    6  = "2026-02-06 22:30:35"
    7  = "2026-02-06 22:30:38"
    8  = "2026-02-06 22:30:40"
    9  = "2026-02-06 22:30:43"
    10 = "2026-02-06 22:30:45"
    11 = "2026-02-06 22:30:47"
    12 = "2026-02-06 22:30:50"
    13 = "2026-02-06 22:30:52"
    14 = "2026-02-06 22:30:54"
    15 = "2026-02-06 22:30:57"
    16 = "2026-02-06 22:30:59"
    17 = "2026-02-06 22:31:02"
    18 = "2026-02-06 22:31:04"
    19 = "2026-02-06 22:31:06"
    20 = "2026-02-06 22:31:09"
    21 = "2026-02-06 22:31:11"
    22 = "2026-02-06 22:31:14"
    23 = "2026-02-06 22:31:16"
    24 = "2026-02-06 22:31:19"
    25 = "2026-02-06 22:31:21"
    26 = "2026-02-06 22:31:24"
    27 = "2026-02-06 22:31:26"
    28 = "2026-02-06 22:31:29"
    29 = "2026-02-06 22:31:31"
    30 = "2026-02-06 22:31:34"
    31 = "2026-02-06 22:31:36"
    32 = "2026-02-06 22:31:39"
    33 = "2026-02-06 22:31:41"
    34 = "2026-02-06 22:31:43"
    35 = "2026-02-06 22:31:46"
    36 = "2026-02-06 22:31:48"
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 5).Value = $updates[$row]
}
